$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) retains text formatting so numeric-looking
# values (e.g. "1.002") are stored as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.744.45"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.907.33"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "239.46"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "0.4943"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("D8").Value = "0.2954"
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("D9").Value = "0.06731"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("D10").Value = "1.907.02"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("D11").Value = "17.01"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").Value = "0.07360"
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").Value = "5.160"
$ws.Range("E13").Value = "  +2.74%  "
$ws.Range("D14").Value = "88.26"
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").Value = "0.6684"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").Value = "30.683.31"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "0.000007917"
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("E18").Value = "  +2.94%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "2.157.57"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "5.299"
$ws.Range("E21").Value = "  +10.53%  "
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "195.03"
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("D24").Value = "6.255"
$ws.Range("E24").Value = "  +2.92%  "
$ws.Range("D25").Value = "9.623"
$ws.Range("E25").Value = "  +2.80%  "
$ws.Range("D26").Value = "162.48"
$ws.Range("E26").Value = "  +3.46%  "
$ws.Range("D27").Value = "18.57"
$ws.Range("E27").Value = "  -1.63%  "
$ws.Range("D28").Value = "1.944"
$ws.Range("E28").Value = "  +2.72%  "
$ws.Range("D29").Value = "1.467"
$ws.Range("E29").Value = "  +4.58%  "
$ws.Range("D30").Value = "4.421"
$ws.Range("E30").Value = "  +3.10%  "
$ws.Range("D31").Value = "0.09133"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("D32").Value = "4.050"
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("D33").Value = "0.05257"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").Value = "0.7405"
$ws.Range("E34").Value = "  +0.41%  "
$ws.Range("D35").Value = "1.109"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("D36").Value = "2.725"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "0.01818"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("D38").Value = "2.715"
$ws.Range("E38").Value = "  +1.41%  "
$ws.Range("D39").Value = "0.9201"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").Value = "2.071"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("D41").Value = "74.15"
$ws.Range("E41").Value = "  +28.84%  "
$ws.Range("D42").Value = "0.4442"
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("E43").Value = "  +1.67%  "
$ws.Range("D44").Value = "5.917"
$ws.Range("E44").Value = "  +3.13%  "
$ws.Range("D45").Value = "0.9986"
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("D46").Value = "0.1386"
$ws.Range("E46").Value = "  +2.99%  "
$ws.Range("D47").Value = "7.568"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").Value = "9.066"
$ws.Range("E48").Value = "  +4.34%  "
$ws.Range("D49").Value = "35.35"
$ws.Range("E49").Value = "  +4.90%  "
$ws.Range("D50").Value = "0.05857"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "0.4001"
$ws.Range("E51").Value = "  +1.57%  "
